$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 922.5625
$ws.Range("I41").Value = 1051.6364
$ws.Range("J41").Value = 638.6
$ws.Range("K41").Value = 1051.6364
$ws.Range("L41").Value = 638.6
$ws.Range("M41").Value = -611.6364000000001
$ws.Range("N41").Value = -1518.6
$ws.Range("H53").Value = 563.1667
$ws.Range("I53").Value = 199
$ws.Range("K53").Value = 199
$ws.Range("M53").Value = 438
$ws.Range("H92").Value = 197.66667
$ws.Range("I92").Value = 216.875
$ws.Range("K92").Value = 216.875
$ws.Range("M92").Value = 1031.125
$ws.Range("H108").Value = 68000
$ws.Range("J108").Value = 68000
$ws.Range("L108").Value = 68000
$ws.Range("N108").Value = -75680
$ws.Range("H137").Value = 43991.875
$ws.Range("I137").Value = 1280.6
$ws.Range("J137").Value = 55231.684
$ws.Range("K137").Value = 3841.8
$ws.Range("L137").Value = 165695.052
$ws.Range("M137").Value = -1291.8
$ws.Range("N137").Value = -170795.052
$ws.Range("H138").Value = 4634.636
$ws.Range("I138").Value = 4634.636
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 13903.908
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -8763.908000000001
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2380.5073
$ws.Range("I32").Value = 1597.6897
$ws.Range("K32").Value = 1597.6897
$ws.Range("M32").Value = -1310.6897
$ws.Range("H45").Value = 2695.0605
$ws.Range("I45").Value = 2603.2307
$ws.Range("J45").Value = 2754.75
$ws.Range("K45").Value = 2603.2307
$ws.Range("L45").Value = 2754.75
$ws.Range("M45").Value = -2226.2307
$ws.Range("N45").Value = -3508.75
$ws.Range("H61").Value = 3655.04
$ws.Range("I61").Value = 2813.4167
$ws.Range("J61").Value = 4431.923
$ws.Range("K61").Value = 2813.4167
$ws.Range("L61").Value = 4431.923
$ws.Range("M61").Value = -2601.4167
$ws.Range("N61").Value = -4855.923
$ws.Range("H102").Value = 1568.0588
$ws.Range("I102").Value = 1510.5333
$ws.Range("K102").Value = 1510.5333
$ws.Range("M102").Value = 111.4666999999999
$ws.Range("H123").Value = 64499
$ws.Range("J123").Value = 64499
$ws.Range("L123").Value = 64499
$ws.Range("N123").Value = -74299
$ws.Range("H132").Value = 2893.3948
$ws.Range("I132").Value = 2737.9
$ws.Range("J132").Value = 3476.5
$ws.Range("K132").Value = 8213.700000000001
$ws.Range("L132").Value = 10429.5
$ws.Range("M132").Value = -5683.700000000001
$ws.Range("N132").Value = -15489.5
$ws.Range("H136").Value = 3655.04
$ws.Range("I136").Value = 2813.4167
$ws.Range("J136").Value = 4431.923
$ws.Range("K136").Value = 8440.250100000001
$ws.Range("L136").Value = 13295.769
$ws.Range("M136").Value = -5890.250100000001
$ws.Range("N136").Value = -18395.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 60000
$ws.Range("J130").Value = 60000
$ws.Range("L130").Value = 60000
$ws.Range("N130").Value = -70040
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("H134").Value = 2975.3333
$ws.Range("I134").Value = 2777.1333
$ws.Range("K134").Value = 8331.3999
$ws.Range("M134").Value = -5796.3999
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 41106.875
$ws.Range("I23").Value = 50000
$ws.Range("K23").Value = 50000
$ws.Range("M23").Value = -49760
$ws.Range("H27").Value = 41106.875
$ws.Range("I27").Value = 50000
$ws.Range("K27").Value = 50000
$ws.Range("M27").Value = -49808
$ws.Range("H31").Value = 1243.326
$ws.Range("I31").Value = 948.4375
$ws.Range("J31").Value = 1400.6
$ws.Range("K31").Value = 948.4375
$ws.Range("L31").Value = 1400.6
$ws.Range("M31").Value = -653.4375
$ws.Range("N31").Value = -1990.6
$ws.Range("H34").Value = 1243.326
$ws.Range("I34").Value = 948.4375
$ws.Range("J34").Value = 1400.6
$ws.Range("K34").Value = 948.4375
$ws.Range("L34").Value = 1400.6
$ws.Range("M34").Value = -746.4375
$ws.Range("N34").Value = -1804.6
$ws.Range("H125").Value = 32111
$ws.Range("J125").Value = 32111
$ws.Range("L125").Value = 32111
$ws.Range("N125").Value = -37031
$ws.Range("H132").Value = 2688.6562
$ws.Range("I132").Value = 2201.9583
$ws.Range("K132").Value = 6605.874899999999
$ws.Range("M132").Value = -4075.874899999999
$ws.Range("H134").Value = 1619.091
$ws.Range("I134").Value = 1057.625
$ws.Range("K134").Value = 3172.875
$ws.Range("M134").Value = -637.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1922.375
$ws.Range("J68").Value = 2048.5117
$ws.Range("L68").Value = 6145.5351
$ws.Range("N68").Value = -7767.5351
$ws.Range("H71").Value = 1922.375
$ws.Range("J71").Value = 2048.5117
$ws.Range("L71").Value = 18436.6053
$ws.Range("N71").Value = -26548.6053
$ws.Range("H107").Value = 1360.2142
$ws.Range("I107").Value = 1292.6666
$ws.Range("J107").Value = 1410.875
$ws.Range("K107").Value = 3877.9998
$ws.Range("L107").Value = 4232.625
$ws.Range("M107").Value = -1957.9998
$ws.Range("N107").Value = -8072.625
$ws.Range("H121").Value = 724.75
$ws.Range("J121").Value = 816.5
$ws.Range("L121").Value = 2449.5
$ws.Range("N121").Value = -5069.5
$ws.Range("H129").Value = 35169.523
$ws.Range("J129").Value = 66606
$ws.Range("L129").Value = 199818
$ws.Range("N129").Value = -209818
$ws.Range("H131").Value = 13534128
$ws.Range("J131").Value = 23056.818
$ws.Range("L131").Value = 69170.454
$ws.Range("N131").Value = -79250.454
$ws.Range("H132").Value = 1044.8334
$ws.Range("I132").Value = 842.8
$ws.Range("J132").Value = 1189.1428
$ws.Range("K132").Value = 7585.2
$ws.Range("L132").Value = 10702.2852
$ws.Range("M132").Value = -5055.2
$ws.Range("N132").Value = -15762.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 16066.667
$ws.Range("J98").Value = 16066.667
$ws.Range("L98").Value = 16066.667
$ws.Range("N98").Value = -22056.667
$ws.Range("H99").Value = 2450
$ws.Range("I99").Value = 2450
$ws.Range("K99").Value = 2450
$ws.Range("M99").Value = -204
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("H126").Value = 1987006.9
$ws.Range("I126").Value = 3475209.8
$ws.Range("J126").Value = 2736.4167
$ws.Range("K126").Value = 10425629.4
$ws.Range("L126").Value = 8209.250100000001
$ws.Range("M126").Value = -10423159.4
$ws.Range("N126").Value = -13149.2501
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 342338
$ws.Range("J43").Value = 342338
$ws.Range("L43").Value = 342338
$ws.Range("N43").Value = -342724
$ws.Range("H122").Value = 7959
$ws.Range("I122").Value = 5662.9287
$ws.Range("J122").Value = 13316.5
$ws.Range("K122").Value = 16988.7861
$ws.Range("L122").Value = 39949.5
$ws.Range("M122").Value = -14538.7861
$ws.Range("N122").Value = -44849.5
$ws.Range("H136").Value = 4937.409
$ws.Range("I136").Value = 3013.6365
$ws.Range("J136").Value = 6861.1816
$ws.Range("K136").Value = 9040.9095
$ws.Range("L136").Value = 20583.5448
$ws.Range("M136").Value = -6490.9095
$ws.Range("N136").Value = -25683.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 1000
$ws.Range("I26").Value = 1000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -707
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("H126").Value = 5018.85
$ws.Range("I126").Value = 4352.4614
$ws.Range("J126").Value = 6256.4287
$ws.Range("K126").Value = 13057.3842
$ws.Range("L126").Value = 18769.2861
$ws.Range("M126").Value = -10587.3842
$ws.Range("N126").Value = -23709.2861
$ws.Range("N26").ClearContents()
$ws.Range("N95").ClearContents()
